# Split the single bibliography run into 14 separate entries, each
# separated by a manual line break (<w:br/>), matching the diff:
# the paragraph previously held one long run of text; now it holds
# one run with 14 <w:t> segments interleaved with <w:br/> elements.
$d = $word.ActiveDocument

$markers = 2..14

foreach ($n in $markers) {
    $needle = "[" + $n + "]"
    $replacement = "^l[" + $n + "]"
    $found = $d.Content.Find.Execute($needle, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $replacement, 2)
    if (-not $found) {
        Write-Host "WARNING: marker not found:" $needle
    }
}
